$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data for the MathNet.Numerics library (Component / Licence / Link)
$ws.Range("A10").Value = "MathNet.Numerics"
$ws.Range("C10").Value = "MIT/X11"
$ws.Range("D10").Value = "https://numerics.mathdotnet.com/License.html"

# The Component (A) and Licence (C) columns elsewhere in the sheet carry a
# boxed border; reproduce a left+right thin border on the new cells.
$ws.Range("A10").Borders.Item(7).LineStyle = 1
$ws.Range("A10").Borders.Item(7).Weight = 2
$ws.Range("A10").Borders.Item(10).LineStyle = 1
$ws.Range("A10").Borders.Item(10).Weight = 2

$ws.Range("C10").Borders.Item(7).LineStyle = 1
$ws.Range("C10").Borders.Item(7).Weight = 2
$ws.Range("C10").Borders.Item(10).LineStyle = 1
$ws.Range("C10").Borders.Item(10).Weight = 2

# Leave the selection where it was when the file was last saved
$ws.Range("E16").Select()
